# Updates from unipept runs to msms summary
# Inserts two new columns ("DB peptides matched to Nitrospina Unipept" at Q,
# and "DNO peptides matched to Nitrospina Unipept" at S) into the "all samples"
# worksheet, shifting the existing Cyano-Unipept / trypsin / DNO columns to the
# right, then fills in the newly measured Unipept-match counts for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the two new columns -----------------------------------------
# Inserting at Q pushes old Q..V -> R..W.
$ws.Columns("Q").Insert()
# Inserting again at (the now-shifted) S pushes old R (currently sitting at S)
# one more step to the right, landing everything on its final column.
$ws.Columns("S").Insert()

# --- 2. New column headers ---------------------------------------------------
$ws.Range("Q1").Value = "DB peptides matched to Nitrospina Unipept"
$ws.Range("S1").Value = "DNO peptides matched to Nitrospina Unipept"

# --- 3. Fill in the new / updated data values --------------------------------

# Row 4 (sample 231)
$ws.Range("Q4").Value = 18
$ws.Range("S4").Value = 0

# Row 6 (sample 233)
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 28
$ws.Range("S6").Value = 1

# Row 14 (sample 243)
$ws.Range("L14").Value = 477
$ws.Range("M14").Value = 405
$ws.Range("P14").Value = 1
$ws.Range("Q14").Value = 7
$ws.Range("S14").Value = 0
$ws.Range("W14").Value = 225

# Row 29 (sample 273)
$ws.Range("L29").Value = 651
$ws.Range("Q29").Value = 0
$ws.Range("S29").Value = 0

# Row 36 (sample 378)
$ws.Range("L36").Value = 330
$ws.Range("M36").Value = 302
$ws.Range("O36").Value = 1207
$ws.Range("P36").Value = 1
$ws.Range("Q36").Value = 0
$ws.Range("R36").Value = 1
$ws.Range("S36").Value = 0
$ws.Range("V36").Value = 551
$ws.Range("W36").Value = 212

# Row 40 (sample 278 / second occurrence)
$ws.Range("L40").Value = 849
$ws.Range("M40").Value = 777
$ws.Range("O40").Value = 9050
$ws.Range("P40").Value = 1
$ws.Range("Q40").Value = 1
$ws.Range("R40").Value = 1
$ws.Range("S40").Value = 1
$ws.Range("V40").Value = 3125
$ws.Range("W40").Value = 522

# --- 4. Restore the view state shown in the saved workbook -------------------
$ws.Range("Q15").Select()
